# Update countries & provincias Spain
# Applies the 26-Jun-2020 data refresh: updates the "last updated" timestamp,
# refreshes case counters for several countries, and re-sorts a handful of
# rows whose ranking changed as a result (which swaps the country name shown
# in those rows together with their case data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the "datos actualizados" timestamp in the title cell (A1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 09:26"

# ---------------------------------------------------------------------
# 2. Countries whose row order changed (ranking swap) together with
#    fresh data for the row that moved up.
# ---------------------------------------------------------------------

# Ucrania overtakes Portugal -> row 37 becomes Ucrania (fresh numbers),
# row 38 becomes Portugal (keeps Portugal's previous numbers).
$ws.Range("A37").Value = "Ucrania"
$ws.Range("B37").Value = 41117
$ws.Range("C37").Value = 1109
$ws.Range("D37").Value = 18299
$ws.Range("E37").Value = 21732
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 19
$ws.Range("H37").Value = 1086

$ws.Range("A38").Value = "Portugal"
$ws.Range("B38").Value = 40415
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 26382
$ws.Range("E38").Value = 12484
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 1549

# Armenia overtakes Nigeria -> row 51 becomes Armenia (fresh numbers),
# row 52 becomes Nigeria (keeps Nigeria's previous numbers).
$ws.Range("A51").Value = "Armenia"
$ws.Range("B51").Value = 23247
$ws.Range("C51").Value = 759
$ws.Range("D51").Value = 12149
$ws.Range("E51").Value = 10688
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 13
$ws.Range("H51").Value = 410

$ws.Range("A52").Value = "Nigeria"
$ws.Range("B52").Value = 22614
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 7822
$ws.Range("E52").Value = 14243
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 549

# Fiyi / Dominica swap places (tied stats, no numeric change).
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

# Groenlandia / Islas Malvinas swap places (tied stats, no numeric change).
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"

# Seychelles overtakes Montserrat -> row 211 becomes Seychelles,
# row 212 becomes Montserrat (their small data sets swap too).
$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

$ws.Range("A212").Value = "Montserrat"
$ws.Range("B212").Value = 11
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 10
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 1

# ---------------------------------------------------------------------
# 3. Plain data refreshes (no reordering involved).
# ---------------------------------------------------------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2504676
$ws.Range("C4").Value = 88
$ws.Range("D4").Value = 1052389
$ws.Range("E4").Value = 1325502
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 126785

# Row 83
$ws.Range("D83").Value = 3147
$ws.Range("E83").Value = 2056
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 7
$ws.Range("H83").Value = 133

# Row 92
$ws.Range("B92").Value = 4127
$ws.Range("C92").Value = 4
$ws.Range("D92").Value = 2663
$ws.Range("E92").Value = 886
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 578

# Row 126
$ws.Range("B126").Value = 1112
$ws.Range("C126").Value = 1
$ws.Range("D126").Value = 932
$ws.Range("E126").Value = 150

# Row 136
$ws.Range("B136").Value = 919
$ws.Range("C136").Value = 2
$ws.Range("D136").Value = 780
$ws.Range("E136").Value = 125
